$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$data  = $wb.Worksheets.Item("Data")
$rpe   = $wb.Worksheets.Item("RPEpUACE")

# --- Data sheet (written first so new shared strings land in the order the
#     target workbook expects) ---------------------------------------------
$data.Range("A1").Value = "Excerpt from Table 6-2:"

$data.Range("B3").Value = 2015
$data.Range("C3").Value = 2016
$data.Range("D3").Value = 2017
$data.Range("E3").Value = 2018
$data.Range("F3").Value = 2019

$data.Range("B4").Value = -791695
$data.Range("C4").Value = -855998
$data.Range("D4").Value = -792046
$data.Range("E4").Value = -824885
$data.Range("F4").Value = -812695

$data.Range("B5").Value = 663
$data.Range("C5").Value = 308
$data.Range("D5").Value = 614
$data.Range("E5").Value = 552
$data.Range("F5").Value = 552

$data.Range("B6").Value = 38
$data.Range("C6").Value = 18
$data.Range("D6").Value = 36
$data.Range("E6").Value = 32
$data.Range("F6").Value = 32

# --- About sheet -------------------------------------------------------
$about.Range("B3").Value = "US EPA"
$about.Range("B5").Value = "Draft Inventory of US Greenhouse Gas Emissions Emissions and Sinks"
$about.Range("B6").Value = "https://www.epa.gov/sites/production/files/2021-02/documents/us-ghg-inventory-2021-main-text.pdf"
$about.Range("B7").Value = "Table 6-3"
$about.Range("B1").Value = "Colorado"

$about.Range("C1").Value = 44515
$about.Range("B4").Value = 2021

# --- Selections (recreate the saved cursor position on each sheet) -----
# Select Data/RPEpUACE first, About last so "About" ends up the active tab
# (matches tabSelected="1" staying on the About sheet in the target file).
[void]$data.Range("F6").Select()
[void]$rpe.Range("B2:B13").Select()
[void]$about.Range("B8").Select()
